# Auto-generated edit script: update cryptos list with new prices/volumes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.875.07"
$ws.Range("E2").Value = "  +7.21%  "

$ws.Range("D3").Value = "3.865.94"
$ws.Range("E3").Value = "  +14.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "425.20"
$ws.Range("E5").Value = "  +11.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.16"
$ws.Range("E6").Value = "  +7.45%  "

$ws.Range("D7").Value = "3.859.43"
$ws.Range("E7").Value = "  +8.86%  "

$ws.Range("E8").Value = "  +7.61%  "

$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.727"
$ws.Range("E10").Value = "  +12.44%  "

$ws.Range("E11").Value = "  +13.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000346"
$ws.Range("E12").Value = "  +19.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.91"
$ws.Range("E13").Value = "  +7.51%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.462.91"
$ws.Range("E14").Value = "  +14.02%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.25"
$ws.Range("E15").Value = "  +13.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.85"
$ws.Range("E16").Value = "  +29.24%  "

$ws.Range("D17").Value = "3.891.47"
$ws.Range("E17").Value = "  +13.74%  "

$ws.Range("E18").Value = "  +1.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.98"
$ws.Range("E19").Value = "  +10.51%  "

$ws.Range("D20").Value = "66.984.52"
$ws.Range("E20").Value = "  +7.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.08"
$ws.Range("E21").Value = "  +8.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "413.67"
$ws.Range("E22").Value = "  +9.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.89"
$ws.Range("E23").Value = "  +12.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.42"
$ws.Range("E24").Value = "  +7.67%  "

$ws.Range("E25").Value = "  +10.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.73"
$ws.Range("E26").Value = "  +15.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.95"
$ws.Range("E27").Value = "  +14.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.24"
$ws.Range("E28").Value = "  +11.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.30"
$ws.Range("E29").Value = "  +2.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.10"
$ws.Range("E30").Value = "  +43.05%  "

$ws.Range("E31").Value = "  +14.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.15"
$ws.Range("E32").Value = "  +13.47%  "

$ws.Range("E33").Value = "  +14.90%  "

$ws.Range("E34").Value = "  +7.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.03"
$ws.Range("E36").Value = "  +7.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.152"
$ws.Range("E37").Value = "  +5.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.80"
$ws.Range("E38").Value = "  -0.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.38"
$ws.Range("E39").Value = "  +37.04%  "

$ws.Range("D40").Value = "0.0₃0760"
$ws.Range("E40").Value = "  +31.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0461"
$ws.Range("E41").Value = "  +9.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.89"
$ws.Range("E42").Value = "  +9.73%  "

$ws.Range("E43").Value = "  +1.30%  "

$ws.Range("B44").Value = "LidoDAOToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.37"
$ws.Range("E44").Value = "  +13.72%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.135"
$ws.Range("E45").Value = "  +4.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.13"
$ws.Range("E46").Value = "  +3.57%  "

$ws.Range("E47").Value = "  +19.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.05"
$ws.Range("E48").Value = "  +7.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.56"
$ws.Range("E49").Value = "  +1.04%  "

$ws.Range("E50").Value = "  +10.22%  "

$ws.Range("E51").Value = "  +5.75%  "
